# Update "想去人数" (F) and "最低票价" (G) figures on both the "展览"
# sheet and the "全部类型" sheet (which mirrors the same rows), to
# reflect refreshed scrape values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 647
    $ws.Range("G2").Value = 68

    $ws.Range("F3").Value = 491

    $ws.Range("F8").Value = 1942

    $ws.Range("F9").Value = 4054

    $ws.Range("F10").Value = 92
}
